$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Beetles" row (row 42) is re-purposed from "boss grasshopper" to "Beetle 1",
# with a new description, as part of designing level 7-1.
$ws.Range("A42").Value = "Beetle 1"
$ws.Range("B36").Value = "Hey Jude."
$ws.Range("B42").Value = "The first born of the Beetle brothers."

# Move the active selection from B46 to B42.
$ws.Range("B42").Select()

# Minimize the workbook window (best-effort; mirrors the authored
# workbookView/@minimized flag in the saved file).
$wb.Windows.Item(1).WindowState = -4140
